$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "Run 50" column (AZ). This shifts the old "Mean" column (BA)
# one position to the left, becoming the new column AZ.
$ws.Columns("AZ").Delete()

# Update header: "Gen" -> "MaxFES"
$ws.Range("A1").Value = "MaxFES"

# Update column A (generation counters -> MaxFES fractions) for rows 2-14
$aValues = New-Object 'object[,]' 13,1
$aValues[0,0] = 0
$aValues[1,0] = 0.001
$aValues[2,0] = 0.01
$aValues[3,0] = 0.1
$aValues[4,0] = 0.2
$aValues[5,0] = 0.3
$aValues[6,0] = 0.4
$aValues[7,0] = 0.5
$aValues[8,0] = 0.6
$aValues[9,0] = 0.7
$aValues[10,0] = 0.8
$aValues[11,0] = 0.9
$aValues[12,0] = 1
$ws.Range("A2:A14").Value = $aValues

# Update the (now shifted) Mean column AZ with the recalculated means
# (Run 50 excluded from the average) for rows 2-14
$meanValues = New-Object 'object[,]' 13,1
$meanValues[0,0] = 1262.25566111
$meanValues[1,0] = 1258.51708948
$meanValues[2,0] = 1016.98382385
$meanValues[3,0] = 376.7103564
$meanValues[4,0] = 133.83568844
$meanValues[5,0] = 37.50267065
$meanValues[6,0] = 8.89236503
$meanValues[7,0] = 2.73508603
$meanValues[8,0] = 1.35486203
$meanValues[9,0] = 1.06656082
$meanValues[10,0] = 0.89059611
$meanValues[11,0] = 0.58498504
$meanValues[12,0] = 0.24470739
$ws.Range("AZ2:AZ14").Value = $meanValues
